# Weekly refresh of the "Vega Modelo de Temuco - Puerro" price series.
# A brand-new record is inserted at the top of the data block (row 105,
# right after the header + the 103 already-"frozen" rows above it),
# pushing every existing record down by one row. The oldest record that
# used to be the last row of the table now becomes a new trailing row.
#
# Net effect on the OOXML: dimension grows from A1:R236 to A1:R237, and a
# new row is inserted at r=105 with the values below (all later rows keep
# their original content, just shifted down by one row index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 105:236 down to 106:237, carrying their values/formats along.
$ws.Range("A105").EntireRow.Insert()

# Populate the newly inserted row 105 with this week's new record.
$ws.Cells.Item(105, 1).Value = 10
$ws.Cells.Item(105, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(105, 3).Value = "La Araucanía"
$ws.Cells.Item(105, 4).Value = 44810
$ws.Cells.Item(105, 5).Value = 9
$ws.Cells.Item(105, 6).Value = 100112005
$ws.Cells.Item(105, 7).Value = "Puerro"
$ws.Cells.Item(105, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 40
$ws.Cells.Item(105, 11).Value = 19000
$ws.Cells.Item(105, 12).Value = 20000
$ws.Cells.Item(105, 13).Value = 19500
$ws.Cells.Item(105, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(105, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(105, 16).Value = 1625
$ws.Cells.Item(105, 17).Value = 12
$ws.Cells.Item(105, 18).Value = "Hortaliza"
